$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.459.41'
$ws.Range("E2").Value = '  +0.42%  '
$ws.Range("D3").Value = '2.108.56'
$ws.Range("E3").Value = '  +0.80%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.005'
$ws.Range("E4").Value = '  +0.39%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '334.81'
$ws.Range("E5").Value = '  +1.85%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.005'
$ws.Range("E6").Value = '  +0.50%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5218'
$ws.Range("E7").Value = '  -0.19%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4524'
$ws.Range("E8").Value = '  +4.16%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '54.18'
$ws.Range("E9").Value = '  +15.87%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.08918'
$ws.Range("E10").Value = '  +1.17%  '
$ws.Range("E11").Value = '  +1.58%  '
$ws.Range("E12").Value = '  -1.40%  '
$ws.Range("D13").Value = '2.104.07'
$ws.Range("E13").Value = '  +0.85%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.825'
$ws.Range("E14").Value = '  +1.32%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '8.018'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '96.76'
$ws.Range("E16").Value = '  +0.50%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001143'
$ws.Range("E17").Value = '  +1.38%  '
$ws.Range("E18").Value = '  +0.50%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06650'
$ws.Range("E19").Value = '  +0.20%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '19.22'
$ws.Range("E20").Value = '  +1.55%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.004'
$ws.Range("E21").Value = '  +0.40%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.329'
$ws.Range("E22").Value = '  -0.16%  '
$ws.Range("D23").Value = '30.549.36'
$ws.Range("E23").Value = '  +0.51%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.42'
$ws.Range("E24").Value = '  +0.31%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.352'
$ws.Range("E25").Value = '  +2.00%  '
$ws.Range("D26").Value = '2.352.41'
$ws.Range("E26").Value = '  +0.92%  '
$ws.Range("E27").Value = '  -1.08%  '
$ws.Range("B28").Value = 'Monero'
$ws.Range("C28").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '162.97'
$ws.Range("E28").Value = '  +0.66%  '
$ws.Range("B29").Value = 'LidoDAOToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.534'
$ws.Range("E29").Value = '  -2.41%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '133.93'
$ws.Range("E30").Value = '  +1.28%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.209'
$ws.Range("E31").Value = '  -0.37%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.1070'
$ws.Range("E32").Value = '  -0.05%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.397'
$ws.Range("E33").Value = '  +3.57%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.642'
$ws.Range("E34").Value = '  -1.84%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.948'
$ws.Range("E35").Value = '  +1.44%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '10.37'
$ws.Range("E36").Value = '  +3.81%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.770'
$ws.Range("E37").Value = '  +5.40%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02584'
$ws.Range("E38").Value = '  -0.09%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06839'
$ws.Range("E39").Value = '  +2.08%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.2304'
$ws.Range("E40").Value = '  +1.85%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '12.74'
$ws.Range("E41").Value = '  +0.67%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.6865'
$ws.Range("E42").Value = '  +0.58%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.247'
$ws.Range("E43").Value = '  -0.16%  '
$ws.Range("B44").Value = 'NEARProtocol'
$ws.Range("C44").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.319'
$ws.Range("E44").Value = '  +4.99%  '
$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '14.01'
$ws.Range("E45").Value = '  -0.14%  '
$ws.Range("B46").Value = 'Decentraland'
$ws.Range("C46").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.6364'
$ws.Range("E46").Value = '  -0.11%  '
$ws.Range("B47").Value = 'PancakeSwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.665'
$ws.Range("E47").Value = '  +1.39%  '
$ws.Range("B48").Value = 'BabyDogeCoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.00000000352'
$ws.Range("E48").Value = '  +23.04%  '
$ws.Range("E49").Value = '  +0.05%  '
$ws.Range("B50").Value = 'WEMIXTOKEN'
$ws.Range("C50").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.206'
$ws.Range("E50").Value = '  +1.19%  '
$ws.Range("B51").Value = 'Aave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '83.13'
$ws.Range("E51").Value = '  +1.47%  '
